# Atualização de bases das ligas, do dia: 06-04-2024 às 15:39
#
# The underlying odds-data refresh re-ordered a handful of already-recorded
# fixtures (rows 236-239) and refreshed the odds for several upcoming
# fixtures (rows 261-265), while one fixture that is no longer tracked
# (the old row 266) was dropped from the sheet entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 236-239: the four fixtures rotate up by one row (column A, the
#     running match index, stays put on each row). Row 239 wraps around and
#     picks up the data that used to sit in row 236, so we must stash that
#     row's original values before anything gets overwritten. ---
$origRow236 = $ws.Range("B236:AC236").Value2

$ws.Range("B236:AC236").Value2 = $ws.Range("B237:AC237").Value2
$ws.Range("B237:AC237").Value2 = $ws.Range("B238:AC238").Value2
$ws.Range("B238:AC238").Value2 = $ws.Range("B239:AC239").Value2
$ws.Range("B239:AC239").Value2 = $origRow236

# --- Rows 261-265: upcoming fixtures shift up by one row, pulling in the
#     row below (which also brings a handful of freshly re-quoted odds).
#     Processing top-to-bottom is safe because each row is read from the
#     row beneath it before that row is touched. ---
$ws.Range("B261:AA261").Value2 = $ws.Range("B262:AA262").Value2
$ws.Range("B262:AA262").Value2 = $ws.Range("B263:AA263").Value2
$ws.Range("B263:AA263").Value2 = $ws.Range("B264:AA264").Value2
$ws.Range("B264:AA264").Value2 = $ws.Range("B265:AA265").Value2
$ws.Range("B265:AA265").Value2 = $ws.Range("B266:AA266").Value2

# A few odds were re-quoted on top of the plain shift above.
$ws.Cells.Item(261, 14).Value2 = 2.1     # N261 oddH
$ws.Cells.Item(261, 16).Value2 = 3.6     # P261 oddA
$ws.Cells.Item(261, 17).Value2 = -0.25   # Q261 Ah
$ws.Cells.Item(261, 18).Value2 = 1.8     # R261 oddAHH
$ws.Cells.Item(261, 19).Value2 = 2.05    # S261 oddAHA
$ws.Cells.Item(261, 20).Value2 = 2.5     # T261 AhOU
$ws.Cells.Item(261, 21).Value2 = 2.05    # U261 oddAHOver
$ws.Cells.Item(261, 22).Value2 = 1.8     # V261 oddAHUnder

$ws.Cells.Item(262, 16).Value2 = 3.4     # P262 oddA
$ws.Cells.Item(262, 21).Value2 = 2.025   # U262 oddAHOver
$ws.Cells.Item(262, 22).Value2 = 1.825   # V262 oddAHUnder

$ws.Cells.Item(263, 16).Value2 = 5       # P263 oddA
$ws.Cells.Item(263, 18).Value2 = 1.9     # R263 oddAHH
$ws.Cells.Item(263, 19).Value2 = 1.95    # S263 oddAHA

# --- The old last fixture (row 266, now duplicated into row 265 above) is
#     removed entirely, shrinking the sheet from 266 to 265 data rows. ---
$ws.Rows.Item(266).Delete()
